$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order match data within rows that were permuted (row index/date/league columns A-E stay fixed) ---
$ws.Cells.Item(25,6).Value = "Polkowice"
$ws.Cells.Item(25,7).Value = 2
$ws.Cells.Item(25,8).Value = "Starowice Dolne"
$ws.Cells.Item(25,9).Value = 2
$ws.Cells.Item(25,10).Value = 1.53
$ws.Cells.Item(25,11).Value = "23/08/2023 11:12"
$ws.Cells.Item(25,12).Value = 1.71
$ws.Cells.Item(25,13).Value = "23/08/2023 16:37"
$ws.Cells.Item(25,14).Value = 4.29
$ws.Cells.Item(25,15).Value = "23/08/2023 11:12"
$ws.Cells.Item(25,16).Value = 4.53
$ws.Cells.Item(25,17).Value = "23/08/2023 16:37"
$ws.Cells.Item(25,18).Value = 4.42
$ws.Cells.Item(25,19).Value = "23/08/2023 11:12"
$ws.Cells.Item(25,20).Value = 3.31
$ws.Cells.Item(25,21).Value = "23/08/2023 16:37"
$ws.Cells.Item(25,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/polkowice-starowice-dolne/txyZrqgm/"
$ws.Cells.Item(26,6).Value = "Kluczbork"
$ws.Cells.Item(26,7).Value = 3
$ws.Cells.Item(26,8).Value = "Warta Gorzow"
$ws.Cells.Item(26,9).Value = 0
$ws.Cells.Item(26,10).Value = 1.85
$ws.Cells.Item(26,11).Value = "22/08/2023 05:12"
$ws.Cells.Item(26,12).Value = 2.36
$ws.Cells.Item(26,13).Value = "23/08/2023 16:57"
$ws.Cells.Item(26,14).Value = 3.38
$ws.Cells.Item(26,15).Value = "22/08/2023 05:12"
$ws.Cells.Item(26,16).Value = 3.09
$ws.Cells.Item(26,17).Value = "23/08/2023 16:56"
$ws.Cells.Item(26,18).Value = 3.18
$ws.Cells.Item(26,19).Value = "22/08/2023 05:12"
$ws.Cells.Item(26,20).Value = 2.8
$ws.Cells.Item(26,21).Value = "23/08/2023 16:57"
$ws.Cells.Item(26,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/kluczbork-warta-gorzow/W2oakoGJ/"
$ws.Cells.Item(27,6).Value = "Gwarek Tarnowskie Gory"
$ws.Cells.Item(27,7).Value = 2
$ws.Cells.Item(27,8).Value = "Goczalkowice Zdroj"
$ws.Cells.Item(27,9).Value = 0
$ws.Cells.Item(27,10).Value = 1.93
$ws.Cells.Item(27,11).Value = "22/08/2023 05:12"
$ws.Cells.Item(27,12).Value = 2.2
$ws.Cells.Item(27,13).Value = "23/08/2023 16:06"
$ws.Cells.Item(27,14).Value = 3.29
$ws.Cells.Item(27,15).Value = "22/08/2023 05:12"
$ws.Cells.Item(27,16).Value = 3.3
$ws.Cells.Item(27,17).Value = "23/08/2023 16:06"
$ws.Cells.Item(27,18).Value = 3.06
$ws.Cells.Item(27,19).Value = "22/08/2023 05:12"
$ws.Cells.Item(27,20).Value = 2.88
$ws.Cells.Item(27,21).Value = "23/08/2023 16:06"
$ws.Cells.Item(27,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/gwarek-tarnowskie-gory-goczalkowice-zdroj/SYWUqPws/"
$ws.Cells.Item(28,6).Value = "Gornik Zabrze II"
$ws.Cells.Item(28,7).Value = 1
$ws.Cells.Item(28,8).Value = "Bytom Odrzanski"
$ws.Cells.Item(28,9).Value = 0
$ws.Cells.Item(28,10).Value = 1.7
$ws.Cells.Item(28,11).Value = "23/08/2023 11:12"
$ws.Cells.Item(28,12).Value = 1.99
$ws.Cells.Item(28,13).Value = "23/08/2023 16:07"
$ws.Cells.Item(28,14).Value = 3.77
$ws.Cells.Item(28,15).Value = "23/08/2023 11:12"
$ws.Cells.Item(28,16).Value = 3.93
$ws.Cells.Item(28,17).Value = "23/08/2023 16:06"
$ws.Cells.Item(28,18).Value = 3.7
$ws.Cells.Item(28,19).Value = "23/08/2023 11:12"
$ws.Cells.Item(28,20).Value = 2.88
$ws.Cells.Item(28,21).Value = "23/08/2023 16:07"
$ws.Cells.Item(28,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/gornik-zabrze-bytom-odrzanski/Agp3l5VP/"
$ws.Cells.Item(29,6).Value = "Zielona Gora"
$ws.Cells.Item(29,7).Value = 0
$ws.Cells.Item(29,8).Value = "Bielsko-Biala"
$ws.Cells.Item(29,9).Value = 2
$ws.Cells.Item(29,10).Value = 3.46
$ws.Cells.Item(29,11).Value = "22/08/2023 05:12"
$ws.Cells.Item(29,12).Value = 3.43
$ws.Cells.Item(29,13).Value = "23/08/2023 16:52"
$ws.Cells.Item(29,14).Value = 3.57
$ws.Cells.Item(29,15).Value = "22/08/2023 05:12"
$ws.Cells.Item(29,16).Value = 3.64
$ws.Cells.Item(29,17).Value = "23/08/2023 16:52"
$ws.Cells.Item(29,18).Value = 1.72
$ws.Cells.Item(29,19).Value = "22/08/2023 05:12"
$ws.Cells.Item(29,20).Value = 1.85
$ws.Cells.Item(29,21).Value = "23/08/2023 16:52"
$ws.Cells.Item(29,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/zielona-gora-rekord-bielsko-biala/AHSku1hC/"
$ws.Cells.Item(47,6).Value = "Sleza Wroclaw"
$ws.Cells.Item(47,7).Value = 3
$ws.Cells.Item(47,8).Value = "Stilon Gorzow"
$ws.Cells.Item(47,9).Value = 1
$ws.Cells.Item(47,10).Value = 1.58
$ws.Cells.Item(47,11).Value = "08/09/2023 04:13"
$ws.Cells.Item(47,12).Value = 1.71
$ws.Cells.Item(47,13).Value = "09/09/2023 15:56"
$ws.Cells.Item(47,14).Value = 3.95
$ws.Cells.Item(47,15).Value = "08/09/2023 04:13"
$ws.Cells.Item(47,16).Value = 3.98
$ws.Cells.Item(47,17).Value = "09/09/2023 15:56"
$ws.Cells.Item(47,18).Value = 3.8
$ws.Cells.Item(47,19).Value = "08/09/2023 04:13"
$ws.Cells.Item(47,20).Value = 3.7
$ws.Cells.Item(47,21).Value = "09/09/2023 15:55"
$ws.Cells.Item(47,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/sleza-wroclaw-stilon-gorzow/jgXnn95d/"
$ws.Cells.Item(48,6).Value = "Warta Gorzow"
$ws.Cells.Item(48,7).Value = 0
$ws.Cells.Item(48,8).Value = "Polkowice"
$ws.Cells.Item(48,9).Value = 0
$ws.Cells.Item(48,10).Value = 2.54
$ws.Cells.Item(48,11).Value = "08/09/2023 04:13"
$ws.Cells.Item(48,12).Value = 2.31
$ws.Cells.Item(48,13).Value = "09/09/2023 15:34"
$ws.Cells.Item(48,14).Value = 3.14
$ws.Cells.Item(48,15).Value = "08/09/2023 04:13"
$ws.Cells.Item(48,16).Value = 3.25
$ws.Cells.Item(48,17).Value = "09/09/2023 14:05"
$ws.Cells.Item(48,18).Value = 2.3
$ws.Cells.Item(48,19).Value = "08/09/2023 04:13"
$ws.Cells.Item(48,20).Value = 2.74
$ws.Cells.Item(48,21).Value = "09/09/2023 15:34"
$ws.Cells.Item(48,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/warta-gorzow-polkowice/8lIo7AsM/"
$ws.Cells.Item(63,6).Value = "Goczalkowice Zdroj"
$ws.Cells.Item(63,7).Value = 2
$ws.Cells.Item(63,8).Value = "Zielona Gora"
$ws.Cells.Item(63,9).Value = 0
$ws.Cells.Item(63,10).Value = 2.04
$ws.Cells.Item(63,11).Value = "22/09/2023 01:12"
$ws.Cells.Item(63,12).Value = 2.34
$ws.Cells.Item(63,13).Value = "23/09/2023 13:57"
$ws.Cells.Item(63,14).Value = 3.18
$ws.Cells.Item(63,15).Value = "22/09/2023 01:12"
$ws.Cells.Item(63,16).Value = 3.04
$ws.Cells.Item(63,17).Value = "23/09/2023 13:51"
$ws.Cells.Item(63,18).Value = 2.91
$ws.Cells.Item(63,19).Value = "22/09/2023 01:12"
$ws.Cells.Item(63,20).Value = 2.87
$ws.Cells.Item(63,21).Value = "23/09/2023 13:57"
$ws.Cells.Item(63,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/goczalkowice-zdroj-zielona-gora/bgi4BPY7/"
$ws.Cells.Item(64,6).Value = "Warta Gorzow"
$ws.Cells.Item(64,7).Value = 1
$ws.Cells.Item(64,8).Value = "Pawlowice"
$ws.Cells.Item(64,9).Value = 3
$ws.Cells.Item(64,10).Value = 2.47
$ws.Cells.Item(64,11).Value = "22/09/2023 01:12"
$ws.Cells.Item(64,12).Value = 2.99
$ws.Cells.Item(64,13).Value = "23/09/2023 13:36"
$ws.Cells.Item(64,14).Value = 3.1
$ws.Cells.Item(64,15).Value = "22/09/2023 01:12"
$ws.Cells.Item(64,16).Value = 3.45
$ws.Cells.Item(64,17).Value = "23/09/2023 13:36"
$ws.Cells.Item(64,18).Value = 2.36
$ws.Cells.Item(64,19).Value = "22/09/2023 01:12"
$ws.Cells.Item(64,20).Value = 2.08
$ws.Cells.Item(64,21).Value = "23/09/2023 13:36"
$ws.Cells.Item(64,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/warta-gorzow-pniowek-pawlowice/QaqjZSQE/"
$ws.Cells.Item(69,6).Value = "Zielona Gora"
$ws.Cells.Item(69,7).Value = 2
$ws.Cells.Item(69,8).Value = "Bytom Odrzanski"
$ws.Cells.Item(69,9).Value = 2
$ws.Cells.Item(69,10).Value = 2.3
$ws.Cells.Item(69,11).Value = "29/09/2023 03:12"
$ws.Cells.Item(69,12).Value = 2.05
$ws.Cells.Item(69,13).Value = "30/09/2023 15:51"
$ws.Cells.Item(69,14).Value = 3.21
$ws.Cells.Item(69,15).Value = "29/09/2023 03:12"
$ws.Cells.Item(69,16).Value = 3.55
$ws.Cells.Item(69,17).Value = "30/09/2023 15:51"
$ws.Cells.Item(69,18).Value = 2.5
$ws.Cells.Item(69,19).Value = "29/09/2023 03:12"
$ws.Cells.Item(69,20).Value = 2.97
$ws.Cells.Item(69,21).Value = "30/09/2023 15:51"
$ws.Cells.Item(69,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/zielona-gora-bytom-odrzanski/nwCsPr37/"
$ws.Cells.Item(70,6).Value = "Carina Gubin"
$ws.Cells.Item(70,7).Value = 1
$ws.Cells.Item(70,8).Value = "Sleza Wroclaw"
$ws.Cells.Item(70,9).Value = 4
$ws.Cells.Item(70,10).Value = 2.43
$ws.Cells.Item(70,11).Value = "29/09/2023 03:12"
$ws.Cells.Item(70,12).Value = 2.47
$ws.Cells.Item(70,13).Value = "30/09/2023 15:58"
$ws.Cells.Item(70,14).Value = 3.32
$ws.Cells.Item(70,15).Value = "29/09/2023 03:12"
$ws.Cells.Item(70,16).Value = 3.68
$ws.Cells.Item(70,17).Value = "30/09/2023 15:58"
$ws.Cells.Item(70,18).Value = 2.31
$ws.Cells.Item(70,19).Value = "29/09/2023 03:12"
$ws.Cells.Item(70,20).Value = 2.34
$ws.Cells.Item(70,21).Value = "30/09/2023 15:58"
$ws.Cells.Item(70,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/carina-gubin-sleza-wroclaw/z13aEn6n/"
$ws.Cells.Item(71,6).Value = "Stilon Gorzow"
$ws.Cells.Item(71,7).Value = 1
$ws.Cells.Item(71,8).Value = "Goczalkowice Zdroj"
$ws.Cells.Item(71,9).Value = 1
$ws.Cells.Item(71,10).Value = 2.53
$ws.Cells.Item(71,11).Value = "29/09/2023 03:12"
$ws.Cells.Item(71,12).Value = 2.9
$ws.Cells.Item(71,13).Value = "30/09/2023 15:58"
$ws.Cells.Item(71,14).Value = 3.16
$ws.Cells.Item(71,15).Value = "29/09/2023 03:12"
$ws.Cells.Item(71,16).Value = 3.45
$ws.Cells.Item(71,17).Value = "30/09/2023 15:58"
$ws.Cells.Item(71,18).Value = 2.3
$ws.Cells.Item(71,19).Value = "29/09/2023 03:12"
$ws.Cells.Item(71,20).Value = 2.12
$ws.Cells.Item(71,21).Value = "30/09/2023 15:58"
$ws.Cells.Item(71,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/stilon-gorzow-goczalkowice-zdroj/QTAoO2ID/"
$ws.Cells.Item(72,6).Value = "Pawlowice"
$ws.Cells.Item(72,7).Value = 3
$ws.Cells.Item(72,8).Value = "Gornik Zabrze II"
$ws.Cells.Item(72,9).Value = 2
$ws.Cells.Item(72,10).Value = 1.91
$ws.Cells.Item(72,11).Value = "29/09/2023 03:12"
$ws.Cells.Item(72,12).Value = 1.72
$ws.Cells.Item(72,13).Value = "30/09/2023 15:07"
$ws.Cells.Item(72,14).Value = 3.45
$ws.Cells.Item(72,15).Value = "29/09/2023 03:12"
$ws.Cells.Item(72,16).Value = 3.9
$ws.Cells.Item(72,17).Value = "30/09/2023 15:07"
$ws.Cells.Item(72,18).Value = 3
$ws.Cells.Item(72,19).Value = "29/09/2023 03:12"
$ws.Cells.Item(72,20).Value = 3.74
$ws.Cells.Item(72,21).Value = "30/09/2023 15:07"
$ws.Cells.Item(72,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/pniowek-pawlowice-gornik-zabrze/vD8ZQ4Yf/"
$ws.Cells.Item(90,6).Value = "Starowice Dolne"
$ws.Cells.Item(90,7).Value = 1
$ws.Cells.Item(90,8).Value = "Sleza Wroclaw"
$ws.Cells.Item(90,9).Value = 2
$ws.Cells.Item(90,10).Value = 3.33
$ws.Cells.Item(90,11).Value = "20/10/2023 00:12"
$ws.Cells.Item(90,12).Value = 4.45
$ws.Cells.Item(90,13).Value = "21/10/2023 12:56"
$ws.Cells.Item(90,14).Value = 3.75
$ws.Cells.Item(90,15).Value = "20/10/2023 00:12"
$ws.Cells.Item(90,16).Value = 4.27
$ws.Cells.Item(90,17).Value = "21/10/2023 12:56"
$ws.Cells.Item(90,18).Value = 1.72
$ws.Cells.Item(90,19).Value = "20/10/2023 00:12"
$ws.Cells.Item(90,20).Value = 1.54
$ws.Cells.Item(90,21).Value = "21/10/2023 12:56"
$ws.Cells.Item(90,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/starowice-dolne-sleza-wroclaw/tYMGw3Zc/"
$ws.Cells.Item(91,6).Value = "Goczalkowice Zdroj"
$ws.Cells.Item(91,7).Value = 3
$ws.Cells.Item(91,8).Value = "Carina Gubin"
$ws.Cells.Item(91,9).Value = 0
$ws.Cells.Item(91,10).Value = 1.76
$ws.Cells.Item(91,11).Value = "20/10/2023 00:12"
$ws.Cells.Item(91,12).Value = 1.83
$ws.Cells.Item(91,13).Value = "21/10/2023 12:03"
$ws.Cells.Item(91,14).Value = 3.41
$ws.Cells.Item(91,15).Value = "20/10/2023 00:12"
$ws.Cells.Item(91,16).Value = 3.47
$ws.Cells.Item(91,17).Value = "21/10/2023 12:03"
$ws.Cells.Item(91,18).Value = 3.46
$ws.Cells.Item(91,19).Value = "20/10/2023 00:12"
$ws.Cells.Item(91,20).Value = 3.67
$ws.Cells.Item(91,21).Value = "21/10/2023 12:03"
$ws.Cells.Item(91,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/goczalkowice-zdroj-carina-gubin/KOLKxNl4/"
$ws.Cells.Item(93,6).Value = "Gwarek Tarnowskie Gory"
$ws.Cells.Item(93,7).Value = 1
$ws.Cells.Item(93,8).Value = "Zielona Gora"
$ws.Cells.Item(93,9).Value = 3
$ws.Cells.Item(93,10).Value = 2.13
$ws.Cells.Item(93,11).Value = "20/10/2023 02:12"
$ws.Cells.Item(93,12).Value = 2.18
$ws.Cells.Item(93,13).Value = "21/10/2023 14:40"
$ws.Cells.Item(93,14).Value = 3.26
$ws.Cells.Item(93,15).Value = "20/10/2023 02:12"
$ws.Cells.Item(93,16).Value = 3.47
$ws.Cells.Item(93,17).Value = "21/10/2023 14:40"
$ws.Cells.Item(93,18).Value = 2.7
$ws.Cells.Item(93,19).Value = "20/10/2023 02:12"
$ws.Cells.Item(93,20).Value = 2.8
$ws.Cells.Item(93,21).Value = "21/10/2023 14:40"
$ws.Cells.Item(93,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/gwarek-tarnowskie-gory-zielona-gora/OGIyZvlT/"
$ws.Cells.Item(94,6).Value = "Polkowice"
$ws.Cells.Item(94,7).Value = 1
$ws.Cells.Item(94,8).Value = "Rakow II"
$ws.Cells.Item(94,9).Value = 2
$ws.Cells.Item(94,10).Value = 1.76
$ws.Cells.Item(94,11).Value = "20/10/2023 02:12"
$ws.Cells.Item(94,12).Value = 1.48
$ws.Cells.Item(94,13).Value = "21/10/2023 14:51"
$ws.Cells.Item(94,14).Value = 3.58
$ws.Cells.Item(94,15).Value = "20/10/2023 02:12"
$ws.Cells.Item(94,16).Value = 4.28
$ws.Cells.Item(94,17).Value = "21/10/2023 14:59"
$ws.Cells.Item(94,18).Value = 3.3
$ws.Cells.Item(94,19).Value = "20/10/2023 02:12"
$ws.Cells.Item(94,20).Value = 5.14
$ws.Cells.Item(94,21).Value = "21/10/2023 14:59"
$ws.Cells.Item(94,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/polkowice-rks-rakow-czestochowa/Uce4Ubdp/"
$ws.Cells.Item(95,6).Value = "Bytom Odrzanski"
$ws.Cells.Item(95,7).Value = 0
$ws.Cells.Item(95,8).Value = "Bielsko-Biala"
$ws.Cells.Item(95,9).Value = 3
$ws.Cells.Item(95,10).Value = 3.43
$ws.Cells.Item(95,11).Value = "20/10/2023 02:12"
$ws.Cells.Item(95,12).Value = 3.07
$ws.Cells.Item(95,13).Value = "21/10/2023 14:58"
$ws.Cells.Item(95,14).Value = 3.57
$ws.Cells.Item(95,15).Value = "20/10/2023 02:12"
$ws.Cells.Item(95,16).Value = 3.42
$ws.Cells.Item(95,17).Value = "21/10/2023 14:58"
$ws.Cells.Item(95,18).Value = 1.74
$ws.Cells.Item(95,19).Value = "20/10/2023 02:12"
$ws.Cells.Item(95,20).Value = 2.05
$ws.Cells.Item(95,21).Value = "21/10/2023 14:58"
$ws.Cells.Item(95,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/bytom-odrzanski-rekord-bielsko-biala/x0GPys4A/"
$ws.Cells.Item(100,6).Value = "Stilon Gorzow"
$ws.Cells.Item(100,7).Value = 2
$ws.Cells.Item(100,8).Value = "Gwarek Tarnowskie Gory"
$ws.Cells.Item(100,9).Value = 1
$ws.Cells.Item(100,10).Value = 2.29
$ws.Cells.Item(100,11).Value = "27/10/2023 02:13"
$ws.Cells.Item(100,12).Value = 2.35
$ws.Cells.Item(100,13).Value = "28/10/2023 13:58"
$ws.Cells.Item(100,14).Value = 3.3
$ws.Cells.Item(100,15).Value = "27/10/2023 02:13"
$ws.Cells.Item(100,16).Value = 3.25
$ws.Cells.Item(100,17).Value = "28/10/2023 13:02"
$ws.Cells.Item(100,18).Value = 2.46
$ws.Cells.Item(100,19).Value = "27/10/2023 02:13"
$ws.Cells.Item(100,20).Value = 2.68
$ws.Cells.Item(100,21).Value = "28/10/2023 13:58"
$ws.Cells.Item(100,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/stilon-gorzow-gwarek-tarnowskie-gory/ADz3k0Y9/"
$ws.Cells.Item(101,6).Value = "Unia Turza Slaska"
$ws.Cells.Item(101,7).Value = 2
$ws.Cells.Item(101,8).Value = "Starowice Dolne"
$ws.Cells.Item(101,9).Value = 0
$ws.Cells.Item(101,10).Value = 1.78
$ws.Cells.Item(101,11).Value = "27/10/2023 02:13"
$ws.Cells.Item(101,12).Value = 1.73
$ws.Cells.Item(101,13).Value = "27/10/2023 16:55"
$ws.Cells.Item(101,14).Value = 3.57
$ws.Cells.Item(101,15).Value = "27/10/2023 02:13"
$ws.Cells.Item(101,16).Value = 3.69
$ws.Cells.Item(101,17).Value = "28/10/2023 12:01"
$ws.Cells.Item(101,18).Value = 3.28
$ws.Cells.Item(101,19).Value = "27/10/2023 02:13"
$ws.Cells.Item(101,20).Value = 3.84
$ws.Cells.Item(101,21).Value = "27/10/2023 16:55"
$ws.Cells.Item(101,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/unia-turza-slaska-starowice-dolne/dYEosxtj/"
$ws.Cells.Item(108,6).Value = "Pawlowice"
$ws.Cells.Item(108,7).Value = 3
$ws.Cells.Item(108,8).Value = "Rakow II"
$ws.Cells.Item(108,9).Value = 0
$ws.Cells.Item(108,10).Value = 1.93
$ws.Cells.Item(108,11).Value = "03/11/2023 02:13"
$ws.Cells.Item(108,12).Value = 1.88
$ws.Cells.Item(108,13).Value = "04/11/2023 13:06"
$ws.Cells.Item(108,14).Value = 3.53
$ws.Cells.Item(108,15).Value = "03/11/2023 02:13"
$ws.Cells.Item(108,16).Value = 3.78
$ws.Cells.Item(108,17).Value = "04/11/2023 13:06"
$ws.Cells.Item(108,18).Value = 2.87
$ws.Cells.Item(108,19).Value = "03/11/2023 02:13"
$ws.Cells.Item(108,20).Value = 3.25
$ws.Cells.Item(108,21).Value = "04/11/2023 13:06"
$ws.Cells.Item(108,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/pniowek-pawlowice-rks-rakow-czestochowa/AgUPs9W0/"
$ws.Cells.Item(109,6).Value = "Polkowice"
$ws.Cells.Item(109,7).Value = 3
$ws.Cells.Item(109,8).Value = "Stilon Gorzow"
$ws.Cells.Item(109,9).Value = 2
$ws.Cells.Item(109,10).Value = 1.59
$ws.Cells.Item(109,11).Value = "03/11/2023 02:13"
$ws.Cells.Item(109,12).Value = 1.57
$ws.Cells.Item(109,13).Value = "04/11/2023 13:52"
$ws.Cells.Item(109,14).Value = 3.85
$ws.Cells.Item(109,15).Value = "03/11/2023 02:13"
$ws.Cells.Item(109,16).Value = 4.07
$ws.Cells.Item(109,17).Value = "04/11/2023 13:52"
$ws.Cells.Item(109,18).Value = 3.83
$ws.Cells.Item(109,19).Value = "03/11/2023 02:13"
$ws.Cells.Item(109,20).Value = 4.49
$ws.Cells.Item(109,21).Value = "04/11/2023 13:52"
$ws.Cells.Item(109,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/polkowice-stilon-gorzow/xtWHqV1l/"
$ws.Cells.Item(110,6).Value = "Slask Wroclaw II"
$ws.Cells.Item(110,7).Value = 2
$ws.Cells.Item(110,8).Value = "Zielona Gora"
$ws.Cells.Item(110,9).Value = 2
$ws.Cells.Item(110,10).Value = 1.53
$ws.Cells.Item(110,11).Value = "03/11/2023 02:13"
$ws.Cells.Item(110,12).Value = 1.67
$ws.Cells.Item(110,13).Value = "04/11/2023 13:09"
$ws.Cells.Item(110,14).Value = 3.96
$ws.Cells.Item(110,15).Value = "03/11/2023 02:13"
$ws.Cells.Item(110,16).Value = 3.99
$ws.Cells.Item(110,17).Value = "04/11/2023 13:09"
$ws.Cells.Item(110,18).Value = 4.1
$ws.Cells.Item(110,19).Value = "03/11/2023 02:13"
$ws.Cells.Item(110,20).Value = 3.9
$ws.Cells.Item(110,21).Value = "04/11/2023 13:09"
$ws.Cells.Item(110,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/slask-wroclaw-zielona-gora/W2VLrkHf/"
$ws.Cells.Item(111,6).Value = "Gwarek Tarnowskie Gory"
$ws.Cells.Item(111,7).Value = 0
$ws.Cells.Item(111,8).Value = "Jelenia Gora"
$ws.Cells.Item(111,9).Value = 1
$ws.Cells.Item(111,10).Value = 2.6
$ws.Cells.Item(111,11).Value = "04/11/2023 12:43"
$ws.Cells.Item(111,12).Value = 2.25
$ws.Cells.Item(111,13).Value = "04/11/2023 13:10"
$ws.Cells.Item(111,14).Value = 3.37
$ws.Cells.Item(111,15).Value = "04/11/2023 12:43"
$ws.Cells.Item(111,16).Value = 3.42
$ws.Cells.Item(111,17).Value = "04/11/2023 13:10"
$ws.Cells.Item(111,18).Value = 2.31
$ws.Cells.Item(111,19).Value = "04/11/2023 12:43"
$ws.Cells.Item(111,20).Value = 2.72
$ws.Cells.Item(111,21).Value = "04/11/2023 13:10"
$ws.Cells.Item(111,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/gwarek-tarnowskie-gory-karkonosze-jelenia-gora/OxSDpBnr/"
$ws.Cells.Item(116,6).Value = "Carina Gubin"
$ws.Cells.Item(116,7).Value = 1
$ws.Cells.Item(116,8).Value = "Gornik Zabrze II"
$ws.Cells.Item(116,9).Value = 2
$ws.Cells.Item(116,10).Value = 2.12
$ws.Cells.Item(116,11).Value = "11/11/2023 02:12"
$ws.Cells.Item(116,12).Value = 2.31
$ws.Cells.Item(116,13).Value = "11/11/2023 12:52"
$ws.Cells.Item(116,14).Value = 3.58
$ws.Cells.Item(116,15).Value = "11/11/2023 02:12"
$ws.Cells.Item(116,16).Value = 3.52
$ws.Cells.Item(116,17).Value = "11/11/2023 12:52"
$ws.Cells.Item(116,18).Value = 2.7
$ws.Cells.Item(116,19).Value = "11/11/2023 02:12"
$ws.Cells.Item(116,20).Value = 2.58
$ws.Cells.Item(116,21).Value = "11/11/2023 12:52"
$ws.Cells.Item(116,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/carina-gubin-gornik-zabrze/tzEqbUGs/"
$ws.Cells.Item(117,6).Value = "Bielsko-Biala"
$ws.Cells.Item(117,7).Value = 3
$ws.Cells.Item(117,8).Value = "Gwarek Tarnowskie Gory"
$ws.Cells.Item(117,9).Value = 1
$ws.Cells.Item(117,10).Value = 1.25
$ws.Cells.Item(117,11).Value = "11/11/2023 01:13"
$ws.Cells.Item(117,12).Value = 1.37
$ws.Cells.Item(117,13).Value = "11/11/2023 12:44"
$ws.Cells.Item(117,14).Value = 5.43
$ws.Cells.Item(117,15).Value = "11/11/2023 01:13"
$ws.Cells.Item(117,16).Value = 5.47
$ws.Cells.Item(117,17).Value = "11/11/2023 12:44"
$ws.Cells.Item(117,18).Value = 7.09
$ws.Cells.Item(117,19).Value = "11/11/2023 01:13"
$ws.Cells.Item(117,20).Value = 5.17
$ws.Cells.Item(117,21).Value = "11/11/2023 12:44"
$ws.Cells.Item(117,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/rekord-bielsko-biala-gwarek-tarnowskie-gory/vRRH2A8Q/"
$ws.Cells.Item(119,6).Value = "Unia Turza Slaska"
$ws.Cells.Item(119,7).Value = 2
$ws.Cells.Item(119,8).Value = "Bytom Odrzanski"
$ws.Cells.Item(119,9).Value = 2
$ws.Cells.Item(119,10).Value = 1.93
$ws.Cells.Item(119,11).Value = "11/11/2023 02:12"
$ws.Cells.Item(119,12).Value = 2.02
$ws.Cells.Item(119,13).Value = "11/11/2023 12:23"
$ws.Cells.Item(119,14).Value = 3.61
$ws.Cells.Item(119,15).Value = "11/11/2023 02:12"
$ws.Cells.Item(119,16).Value = 3.41
$ws.Cells.Item(119,17).Value = "11/11/2023 12:23"
$ws.Cells.Item(119,18).Value = 3.05
$ws.Cells.Item(119,19).Value = "11/11/2023 02:12"
$ws.Cells.Item(119,20).Value = 3.15
$ws.Cells.Item(119,21).Value = "11/11/2023 12:23"
$ws.Cells.Item(119,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/unia-turza-slaska-bytom-odrzanski/l0Hid8of/"
$ws.Cells.Item(120,6).Value = "Zielona Gora"
$ws.Cells.Item(120,7).Value = 1
$ws.Cells.Item(120,8).Value = "Pawlowice"
$ws.Cells.Item(120,9).Value = 1
$ws.Cells.Item(120,10).Value = 2.06
$ws.Cells.Item(120,11).Value = "11/11/2023 01:13"
$ws.Cells.Item(120,12).Value = 1.84
$ws.Cells.Item(120,13).Value = "11/11/2023 08:02"
$ws.Cells.Item(120,14).Value = 3.48
$ws.Cells.Item(120,15).Value = "11/11/2023 01:13"
$ws.Cells.Item(120,16).Value = 3.59
$ws.Cells.Item(120,17).Value = "11/11/2023 11:04"
$ws.Cells.Item(120,18).Value = 2.87
$ws.Cells.Item(120,19).Value = "11/11/2023 01:13"
$ws.Cells.Item(120,20).Value = 3.47
$ws.Cells.Item(120,21).Value = "11/11/2023 08:02"
$ws.Cells.Item(120,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/zielona-gora-pniowek-pawlowice/UepRum1D/"

# --- Append new match rows 132-138 ---
$ws.Cells.Item(132,1).Value = 131
$ws.Cells.Item(132,2).Value = "poland"
$ws.Cells.Item(132,3).Value = "iii-liga-group-iii"
$ws.Cells.Item(132,4).Value = "2023-2024"
$ws.Cells.Item(132,5).Value = 45255.5
$ws.Cells.Item(132,6).Value = "Bytom Odrzanski"
$ws.Cells.Item(132,7).Value = 1
$ws.Cells.Item(132,8).Value = "Goczalkowice Zdroj"
$ws.Cells.Item(132,9).Value = 1
$ws.Cells.Item(132,10).Value = 2.66
$ws.Cells.Item(132,11).Value = "25/11/2023 01:13"
$ws.Cells.Item(132,12).Value = 2.68
$ws.Cells.Item(132,13).Value = "25/11/2023 11:58"
$ws.Cells.Item(132,14).Value = 3.29
$ws.Cells.Item(132,15).Value = "25/11/2023 01:13"
$ws.Cells.Item(132,16).Value = 3.61
$ws.Cells.Item(132,17).Value = "25/11/2023 11:58"
$ws.Cells.Item(132,18).Value = 2.27
$ws.Cells.Item(132,19).Value = "25/11/2023 01:13"
$ws.Cells.Item(132,20).Value = 2.19
$ws.Cells.Item(132,21).Value = "25/11/2023 11:58"
$ws.Cells.Item(132,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/bytom-odrzanski-goczalkowice-zdroj/bD5CM4aa/"

$ws.Cells.Item(133,1).Value = 132
$ws.Cells.Item(133,2).Value = "poland"
$ws.Cells.Item(133,3).Value = "iii-liga-group-iii"
$ws.Cells.Item(133,4).Value = "2023-2024"
$ws.Cells.Item(133,5).Value = 45255.5
$ws.Cells.Item(133,6).Value = "Gornik Zabrze II"
$ws.Cells.Item(133,7).Value = 4
$ws.Cells.Item(133,8).Value = "Unia Turza Slaska"
$ws.Cells.Item(133,9).Value = 0
$ws.Cells.Item(133,10).Value = 2.06
$ws.Cells.Item(133,11).Value = "25/11/2023 01:13"
$ws.Cells.Item(133,12).Value = 2.08
$ws.Cells.Item(133,13).Value = "25/11/2023 01:34"
$ws.Cells.Item(133,14).Value = 3.65
$ws.Cells.Item(133,15).Value = "25/11/2023 01:13"
$ws.Cells.Item(133,16).Value = 3.59
$ws.Cells.Item(133,17).Value = "25/11/2023 10:01"
$ws.Cells.Item(133,18).Value = 2.83
$ws.Cells.Item(133,19).Value = "25/11/2023 01:13"
$ws.Cells.Item(133,20).Value = 2.85
$ws.Cells.Item(133,21).Value = "25/11/2023 03:11"
$ws.Cells.Item(133,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/gornik-zabrze-unia-turza-slaska/Qa74OQUn/"

$ws.Cells.Item(134,1).Value = 133
$ws.Cells.Item(134,2).Value = "poland"
$ws.Cells.Item(134,3).Value = "iii-liga-group-iii"
$ws.Cells.Item(134,4).Value = "2023-2024"
$ws.Cells.Item(134,5).Value = 45255.54166666666
$ws.Cells.Item(134,6).Value = "Polkowice"
$ws.Cells.Item(134,7).Value = 4
$ws.Cells.Item(134,8).Value = "Carina Gubin"
$ws.Cells.Item(134,9).Value = 1
$ws.Cells.Item(134,10).Value = 1.57
$ws.Cells.Item(134,11).Value = "25/11/2023 02:12"
$ws.Cells.Item(134,12).Value = 1.57
$ws.Cells.Item(134,13).Value = "25/11/2023 12:59"
$ws.Cells.Item(134,14).Value = 4.08
$ws.Cells.Item(134,15).Value = "25/11/2023 02:12"
$ws.Cells.Item(134,16).Value = 4.15
$ws.Cells.Item(134,17).Value = "25/11/2023 12:59"
$ws.Cells.Item(134,18).Value = 4.12
$ws.Cells.Item(134,19).Value = "25/11/2023 02:12"
$ws.Cells.Item(134,20).Value = 4.36
$ws.Cells.Item(134,21).Value = "25/11/2023 12:59"
$ws.Cells.Item(134,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/polkowice-carina-gubin/6FQXTn7P/"

$ws.Cells.Item(135,1).Value = 134
$ws.Cells.Item(135,2).Value = "poland"
$ws.Cells.Item(135,3).Value = "iii-liga-group-iii"
$ws.Cells.Item(135,4).Value = "2023-2024"
$ws.Cells.Item(135,5).Value = 45255.54166666666
$ws.Cells.Item(135,6).Value = "Warta Gorzow"
$ws.Cells.Item(135,7).Value = 4
$ws.Cells.Item(135,8).Value = "Starowice Dolne"
$ws.Cells.Item(135,9).Value = 4
$ws.Cells.Item(135,10).Value = 1.71
$ws.Cells.Item(135,11).Value = "25/11/2023 02:12"
$ws.Cells.Item(135,12).Value = 1.68
$ws.Cells.Item(135,13).Value = "25/11/2023 05:05"
$ws.Cells.Item(135,14).Value = 3.69
$ws.Cells.Item(135,15).Value = "25/11/2023 02:12"
$ws.Cells.Item(135,16).Value = 3.8
$ws.Cells.Item(135,17).Value = "25/11/2023 11:04"
$ws.Cells.Item(135,18).Value = 3.74
$ws.Cells.Item(135,19).Value = "25/11/2023 02:12"
$ws.Cells.Item(135,20).Value = 3.98
$ws.Cells.Item(135,21).Value = "25/11/2023 05:05"
$ws.Cells.Item(135,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/warta-gorzow-starowice-dolne/4j68Nppg/"

$ws.Cells.Item(136,1).Value = 135
$ws.Cells.Item(136,2).Value = "poland"
$ws.Cells.Item(136,3).Value = "iii-liga-group-iii"
$ws.Cells.Item(136,4).Value = "2023-2024"
$ws.Cells.Item(136,5).Value = 45255.5625
$ws.Cells.Item(136,6).Value = "Pawlowice"
$ws.Cells.Item(136,7).Value = 2
$ws.Cells.Item(136,8).Value = "Jelenia Gora"
$ws.Cells.Item(136,9).Value = 3
$ws.Cells.Item(136,10).Value = 2.08
$ws.Cells.Item(136,11).Value = "25/11/2023 02:43"
$ws.Cells.Item(136,12).Value = 2.39
$ws.Cells.Item(136,13).Value = "25/11/2023 08:19"
$ws.Cells.Item(136,14).Value = 3.59
$ws.Cells.Item(136,15).Value = "25/11/2023 02:43"
$ws.Cells.Item(136,16).Value = 3.48
$ws.Cells.Item(136,17).Value = "25/11/2023 11:31"
$ws.Cells.Item(136,18).Value = 2.76
$ws.Cells.Item(136,19).Value = "25/11/2023 02:43"
$ws.Cells.Item(136,20).Value = 2.47
$ws.Cells.Item(136,21).Value = "25/11/2023 08:19"
$ws.Cells.Item(136,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/pniowek-pawlowice-karkonosze-jelenia-gora/v9Yo6qFI/"

$ws.Cells.Item(137,1).Value = 136
$ws.Cells.Item(137,2).Value = "poland"
$ws.Cells.Item(137,3).Value = "iii-liga-group-iii"
$ws.Cells.Item(137,4).Value = "2023-2024"
$ws.Cells.Item(137,5).Value = 45256.45833333334
$ws.Cells.Item(137,6).Value = "Rakow II"
$ws.Cells.Item(137,7).Value = 1
$ws.Cells.Item(137,8).Value = "Stilon Gorzow"
$ws.Cells.Item(137,9).Value = 3
$ws.Cells.Item(137,10).Value = 2.01
$ws.Cells.Item(137,11).Value = "26/11/2023 00:13"
$ws.Cells.Item(137,12).Value = 2.61
$ws.Cells.Item(137,13).Value = "26/11/2023 10:46"
$ws.Cells.Item(137,14).Value = 3.66
$ws.Cells.Item(137,15).Value = "26/11/2023 00:13"
$ws.Cells.Item(137,16).Value = 3.53
$ws.Cells.Item(137,17).Value = "26/11/2023 10:46"
$ws.Cells.Item(137,18).Value = 2.86
$ws.Cells.Item(137,19).Value = "26/11/2023 00:13"
$ws.Cells.Item(137,20).Value = 2.29
$ws.Cells.Item(137,21).Value = "26/11/2023 10:46"
$ws.Cells.Item(137,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/rks-rakow-czestochowa-stilon-gorzow/dpqs7PaC/"

$ws.Cells.Item(138,1).Value = 137
$ws.Cells.Item(138,2).Value = "poland"
$ws.Cells.Item(138,3).Value = "iii-liga-group-iii"
$ws.Cells.Item(138,4).Value = "2023-2024"
$ws.Cells.Item(138,5).Value = 45256.55208333334
$ws.Cells.Item(138,6).Value = "Slask Wroclaw II"
$ws.Cells.Item(138,7).Value = 3
$ws.Cells.Item(138,8).Value = "Bielsko-Biala"
$ws.Cells.Item(138,9).Value = 1
$ws.Cells.Item(138,10).Value = 2.31
$ws.Cells.Item(138,11).Value = "26/11/2023 02:42"
$ws.Cells.Item(138,12).Value = 2.22
$ws.Cells.Item(138,13).Value = "26/11/2023 13:03"
$ws.Cells.Item(138,14).Value = 3.43
$ws.Cells.Item(138,15).Value = "26/11/2023 02:42"
$ws.Cells.Item(138,16).Value = 3.55
$ws.Cells.Item(138,17).Value = "26/11/2023 13:03"
$ws.Cells.Item(138,18).Value = 2.52
$ws.Cells.Item(138,19).Value = "26/11/2023 02:42"
$ws.Cells.Item(138,20).Value = 2.68
$ws.Cells.Item(138,21).Value = "26/11/2023 13:03"
$ws.Cells.Item(138,22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/slask-wroclaw-rekord-bielsko-biala/hCxj53UO/"

